$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 161, shifting existing rows 161-252 down to 162-253.
$ws.Rows(161).Insert()

# Populate the newly inserted row 161 with the new data record.
$ws.Range("A161").Value = 8
$ws.Range("B161").Value = "Terminal La Palmera de La Serena"
$ws.Range("C161").Value = "Coquimbo"
$ws.Range("D161").Value = 44488
$ws.Range("E161").Value = 4
$ws.Range("F161").Value = 100114001
$ws.Range("G161").Value = "Papa"
$ws.Range("H161").Value = "Cardinal"
$ws.Range("I161").Value = "1a (cosecha)"
$ws.Range("J161").Value = 2600
$ws.Range("K161").Value = 12000
$ws.Range("L161").Value = 13000
$ws.Range("M161").Value = 12500
$ws.Range("N161").Value = "`$/saco 25 kilos"
$ws.Range("O161").Value = "Provincia del Elquí"
$ws.Range("P161").Value = 500
$ws.Range("Q161").Value = 25
$ws.Range("R161").Value = "Hortaliza"
